# Reorder the "Recorded By" (column G) entries on the active sheet so that
# any "System"/"system" tokens are moved to the front of the comma-separated
# list, while the remaining tokens (e.g. email addresses) keep their
# original relative order and follow after.
#
# Example: "backup@backdoor.com, System, system" -> "System, system, backup@backdoor.com"
#          "dnasr281@gmail.com, System"           -> "System, dnasr281@gmail.com"
#
# Rows whose list does NOT end with a "system" token (e.g. already starts
# with "System", or has no "System" token at all) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($null -eq $value -or $value -eq "") {
        continue
    }

    $parts = $value -split ",\s*" | ForEach-Object { $_.Trim() }

    if ($parts.Count -le 1) {
        continue
    }

    $lastPart = $parts[$parts.Count - 1]
    if ($lastPart.ToLower() -ne "system") {
        continue
    }

    $systemParts = @()
    $otherParts = @()
    foreach ($part in $parts) {
        if ($part.ToLower() -eq "system") {
            $systemParts += $part
        } else {
            $otherParts += $part
        }
    }

    $newValue = ($systemParts + $otherParts) -join ", "
    $cell.Value = $newValue
}
